# Updated cryptos list on Mon Jul 22 10:40:52 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $c = $ws.Range($addr)
    # Force text storage (these display-formatted numbers, e.g. "1.00" or
    # "67.451.86", must stay text, not be re-interpreted as numeric values)
    $c.NumberFormat = "@"
    $c.Value = $val
    # Drop the now-unneeded explicit "Text" number format so the cell keeps
    # its original (default) style, matching how the source data was built.
    $c.ClearFormats()
}

Set-TextValue $ws "D2" "67.451.86"
Set-TextValue $ws "E2" "  +0.78%  "

Set-TextValue $ws "D3" "3.494.06"
Set-TextValue $ws "E3" "  -0.14%  "

Set-TextValue $ws "D4" "1.00"
Set-TextValue $ws "E4" "  +0.00%  "

Set-TextValue $ws "D5" "597.23"
Set-TextValue $ws "E5" "  +0.47%  "

Set-TextValue $ws "D6" "179.50"
Set-TextValue $ws "E6" "  +3.94%  "

Set-TextValue $ws "D7" "1.00"
Set-TextValue $ws "E7" "  +0.02%  "

Set-TextValue $ws "D8" "0.604"
Set-TextValue $ws "E8" "  +2.25%  "

Set-TextValue $ws "D9" "3.495.38"
Set-TextValue $ws "E9" "  -0.07%  "

Set-TextValue $ws "E10" "  +4.89%  "

Set-TextValue $ws "E11" "  -2.10%  "

Set-TextValue $ws "D12" "0.437"
Set-TextValue $ws "E12" "  +1.14%  "

Set-TextValue $ws "D13" "4.086.82"
Set-TextValue $ws "E13" "  -0.36%  "

Set-TextValue $ws "D14" "32.32"
Set-TextValue $ws "E14" "  +10.48%  "

Set-TextValue $ws "D15" "0.135"
Set-TextValue $ws "E15" "  +0.81%  "

Set-TextValue $ws "D16" "67.414.31"
Set-TextValue $ws "E16" "  +0.71%  "

Set-TextValue $ws "E17" "  -0.67%  "

Set-TextValue $ws "D18" "3.476.83"
Set-TextValue $ws "E18" "  -0.56%  "

Set-TextValue $ws "D20" "14.33"
Set-TextValue $ws "E20" "  +0.30%  "

Set-TextValue $ws "D21" "389.98"
Set-TextValue $ws "E21" "  -1.36%  "

Set-TextValue $ws "D22" "7.94"
Set-TextValue $ws "E22" "  -0.21%  "

Set-TextValue $ws "D23" "74.10"
Set-TextValue $ws "E23" "  +0.88%  "

Set-TextValue $ws "E24" "  +1.17%  "

Set-TextValue $ws "E25" "  +0.01%  "

Set-TextValue $ws "D26" "5.73"
Set-TextValue $ws "E26" "  +0.80%  "

Set-TextValue $ws "D27" "0.0000122"
Set-TextValue $ws "E27" "  +0.20%  "

Set-TextValue $ws "D28" "10.40"
Set-TextValue $ws "E28" "  +1.55%  "

Set-TextValue $ws "D29" "0.176"
Set-TextValue $ws "E29" "  -2.88%  "

Set-TextValue $ws "D30" "0.999"
Set-TextValue $ws "E30" "  +0.14%  "

Set-TextValue $ws "D31" "6.25"
Set-TextValue $ws "E31" "  +0.98%  "

Set-TextValue $ws "E32" "  +0.12%  "

Set-TextValue $ws "E33" "  +0.41%  "

Set-TextValue $ws "E34" "  -0.69%  "

Set-TextValue $ws "E35" "  +0.38%  "

Set-TextValue $ws "E36" "  +0.02%  "

Set-TextValue $ws "E37" "  -0.89%  "

Set-TextValue $ws "D38" "164.99"
Set-TextValue $ws "E38" "  +1.70%  "

Set-TextValue $ws "D39" "0.871"
Set-TextValue $ws "E39" "  -1.17%  "

Set-TextValue $ws "E40" "  +9.95%  "

Set-TextValue $ws "E41" "  -1.23%  "

Set-TextValue $ws "E42" "  -1.67%  "

Set-TextValue $ws "E43" "  +0.25%  "

Set-TextValue $ws "D44" "2.849.17"
Set-TextValue $ws "E44" "  +0.30%  "

Set-TextValue $ws "D45" "27.06"
Set-TextValue $ws "E45" "  -0.19%  "

Set-TextValue $ws "D46" "26.28"
Set-TextValue $ws "E46" "  +0.38%  "

Set-TextValue $ws "D47" "0.0724"
Set-TextValue $ws "E47" "  -2.09%  "

Set-TextValue $ws "D48" "41.72"
Set-TextValue $ws "E48" "  -2.40%  "

Set-TextValue $ws "E49" "  -0.72%  "

Set-TextValue $ws "D50" "336.48"
Set-TextValue $ws "E50" "  -0.13%  "

Set-TextValue $ws "D51" "1.06"
Set-TextValue $ws "E51" "  -1.45%  "
